# Added Test cases in BWP Bootstrap
# The edit removes two "blank value / verify error" test-case rows
# (rows 8 & 9: "UDF 5 required (dropdown)..." and "UDF 3, 8 required (textbox)...")
# from the EmulatorData sheet. The remaining rows below shift up, and the
# ID column (B) is renumbered sequentially.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmulatorData")

# Delete the two rows that were removed. Deleting row 8 twice removes
# both the original row 8 and row 9, shifting everything below up by two.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# After the shift:
#  - old row 10 (ID 10) is now row 8  -> renumber to 7
#  - old row 11 (ID 11) is now row 9  -> renumber to 8
#  - old row 12 (ID 12) is now row 10 -> renumber to 9
#  - old row 13 (ID 13) is now row 11 -> renumber to 10
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(11, 2).Value = 10

# Update the selected/active cell to match the post-edit selection
$ws.Range("B11").Select()
